# Updates the cryptos price table (Price + Volume(1h) columns, and one
# coin's Name/Link/Price/Volume in row 51) to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold text that often *looks*
# numeric (e.g. "587.31", "  +1.02%  "). Force the range to Text format
# before writing so Excel doesn't silently coerce it to a number, then
# restore the default style so no stray formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '69.374.27'
$ws.Range("E2").Value = '  +2.11%  '
$ws.Range("D3").Value = '3.390.80'
$ws.Range("E3").Value = '  +1.84%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '587.31'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").Value = '179.56'
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.596'
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("E9").Value = '  +6.18%  '
$ws.Range("D10").Value = '0.590'
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").Value = '48.49'
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("D12").Value = '0.0000282'
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("D13").Value = '677.85'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '8.61'
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("D15").Value = '3.933.44'
$ws.Range("E15").Value = '  +1.62%  '
$ws.Range("D16").Value = '69.373.81'
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").Value = '3.386.85'
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").Value = '17.64'
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").Value = '11.24'
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D21").Value = '0.903'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = '17.10'
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").Value = '103.30'
$ws.Range("E24").Value = '  +4.06%  '
$ws.Range("D25").Value = '3.92'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").Value = '9.71'
$ws.Range("E27").Value = '  +1.56%  '
$ws.Range("D28").Value = '34.09'
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("D29").Value = '8.70'
$ws.Range("E29").Value = '  +1.45%  '
$ws.Range("D30").Value = '6.98'
$ws.Range("E30").Value = '  -1.56%  '
$ws.Range("D31").Value = '11.16'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("D32").Value = '555.27'
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("D33").Value = '3.59'
$ws.Range("E33").Value = '  +6.69%  '
$ws.Range("E34").Value = '  +0.69%  '
$ws.Range("D35").Value = '58.07'
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").Value = '3.689.83'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("E38").Value = '  +5.36%  '
$ws.Range("D39").Value = '35.06'
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("D40").Value = '3.26'
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("D41").Value = '2.68'
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("D42").Value = '0.0₃0697'
$ws.Range("E42").Value = '  +3.29%  '
$ws.Range("D43").Value = '0.338'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("D44").Value = '0.0422'
$ws.Range("E44").Value = '  +3.75%  '
$ws.Range("D45").Value = '3.30'
$ws.Range("E45").Value = '  -1.62%  '
$ws.Range("D46").Value = '2.67'
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("D47").Value = '0.130'
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("E48").Value = '  +5.16%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '132.32'
$ws.Range("E50").Value = '  +1.56%  '
$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").Value = '2.58'
$ws.Range("E51").Value = '  +2.15%  '

$dataRange.Style = "Normal"
